$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "217.25") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.176.98'
$ws.Range("E2").Value = '  -6.48%  '
$ws.Range("D3").Value = '1.666.84'
$ws.Range("E3").Value = '  -4.28%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '217.25'
$ws.Range("E5").Value = '  -4.17%  '
$ws.Range("D6").Value = '0.5097'
$ws.Range("E6").Value = '  -12.17%  '
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").Value = '0.2648'
$ws.Range("E8").Value = '  -3.09%  '
$ws.Range("D9").Value = '0.06343'
$ws.Range("E9").Value = '  -4.27%  '
$ws.Range("D10").Value = '21.54'
$ws.Range("E10").Value = '  -7.68%  '
$ws.Range("D11").Value = '0.07361'
$ws.Range("E11").Value = '  -2.67%  '
$ws.Range("D12").Value = '1.670.76'
$ws.Range("E12").Value = '  -4.25%  '
$ws.Range("D13").Value = '4.541'
$ws.Range("E13").Value = '  -3.51%  '
$ws.Range("D14").Value = '0.5770'
$ws.Range("E14").Value = '  -4.44%  '
$ws.Range("D15").Value = '1.893.68'
$ws.Range("E15").Value = '  -4.29%  '
$ws.Range("D16").Value = '0.000008515'
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("D17").Value = '64.89'
$ws.Range("E17").Value = '  -13.09%  '
$ws.Range("D18").Value = '26.239.42'
$ws.Range("E18").Value = '  -6.25%  '
$ws.Range("D19").Value = '4.928'
$ws.Range("E19").Value = '  -7.56%  '
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").Value = '10.83'
$ws.Range("E21").Value = '  -4.18%  '
$ws.Range("D22").Value = '188.71'
$ws.Range("E22").Value = '  -8.29%  '
$ws.Range("D23").Value = '6.176'
$ws.Range("E23").Value = '  -6.87%  '
$ws.Range("D24").Value = '1.006'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").Value = '143.01'
$ws.Range("E25").Value = '  -4.83%  '
$ws.Range("D26").Value = '7.639'
$ws.Range("E26").Value = '  -5.81%  '
$ws.Range("D27").Value = '0.1171'
$ws.Range("E27").Value = '  -5.22%  '
$ws.Range("D28").Value = '15.69'
$ws.Range("E28").Value = '  -2.86%  '
$ws.Range("D29").Value = '0.05798'
$ws.Range("E29").Value = '  -5.87%  '
$ws.Range("D30").Value = '1.284'
$ws.Range("E30").Value = '  -6.84%  '
$ws.Range("D31").Value = '1.321'
$ws.Range("E31").Value = '  -5.24%  '
$ws.Range("D32").Value = '3.520'
$ws.Range("E32").Value = '  -5.57%  '
$ws.Range("D33").Value = '3.500'
$ws.Range("E33").Value = '  -6.52%  '
$ws.Range("D34").Value = '1.646'
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").Value = '1.005'
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("D36").Value = '0.5984'
$ws.Range("E36").Value = '  -6.72%  '
$ws.Range("D37").Value = '2.355'
$ws.Range("E37").Value = '  -2.58%  '
$ws.Range("D38").Value = '2.637'
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("D39").Value = '0.01607'
$ws.Range("E39").Value = '  -3.73%  '
$ws.Range("D40").Value = '5.993'
$ws.Range("E40").Value = '  -3.10%  '
$ws.Range("D41").Value = '1.081.69'
$ws.Range("E41").Value = '  -4.49%  '
$ws.Range("D42").Value = '0.8578'
$ws.Range("E42").Value = '  -2.28%  '
$ws.Range("D43").Value = '1.008'
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").Value = '99.85'
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("D45").Value = '1.817.00'
$ws.Range("E45").Value = '  -3.96%  '
$ws.Range("D46").Value = '0.00000000110'
$ws.Range("E46").Value = '  +1.44%  '
$ws.Range("D47").Value = '55.68'
$ws.Range("E47").Value = '  -6.39%  '
$ws.Range("D48").Value = '1.006'
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("D49").Value = '8.065'
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("D50").Value = '0.4296'
$ws.Range("E50").Value = '  -2.78%  '
$ws.Range("D51").Value = '0.05179'
$ws.Range("E51").Value = '  -3.77%  '

# Restore the default cell style (no explicit style index), matching
# the original workbook formatting for these cells.
$ws.Range("D2:E51").Style = "Normal"
